$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First set the label columns (A-D) for all rows, choosing an assignment
#     order so that newly-introduced strings ("MuSCs" then "Resolving-Mac")
#     first appear in that relative order. ---
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pomc"
$ws.Cells.Item(2, 3).Value = "Oprd1"
$ws.Cells.Item(2, 4).Value = "FAPs"

$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 2).Value = "Pomc"
$ws.Cells.Item(4, 3).Value = "Oprd1"
$ws.Cells.Item(4, 4).Value = "FAPs"

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pomc"
$ws.Cells.Item(3, 3).Value = "Oprd1"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"

$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Pomc"
$ws.Cells.Item(5, 3).Value = "Oprd1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"

# --- Row 2: recomputed numeric metrics (ECs -> Pomc/Oprd1 -> FAPs) ---
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.5131323333333334
$ws.Cells.Item(2, 8).Value = 1.539397
$ws.Cells.Item(2, 9).Value = 0.9119558630037493
$ws.Cells.Item(2, 10).Value = 0.9119558630037494
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.02165966666666666
$ws.Cells.Item(2, 14).Value = 0.064979
$ws.Cells.Item(2, 15).Value = 0.7290280598220596
$ws.Cells.Item(2, 16).Value = 0.7290280598220598
$ws.Cells.Item(2, 17).Value = 0.01111427529588889
$ws.Cells.Item(2, 18).Value = 0.100028477663
$ws.Cells.Item(2, 19).Value = 0.6648414134489754
$ws.Cells.Item(2, 20).Value = 0.6648414134489755

# --- Row 3 (new): ECs -> Pomc/Oprd1 -> Resolving-Mac ---
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.5131323333333334
$ws.Cells.Item(3, 8).Value = 1.539397
$ws.Cells.Item(3, 9).Value = 0.9119558630037493
$ws.Cells.Item(3, 10).Value = 0.9119558630037494
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.008050666666666666
$ws.Cells.Item(3, 14).Value = 0.024152
$ws.Cells.Item(3, 15).Value = 0.2709719401779404
$ws.Cells.Item(3, 16).Value = 0.2709719401779404
$ws.Cells.Item(3, 17).Value = 0.004131057371555556
$ws.Cells.Item(3, 18).Value = 0.037179516344
$ws.Cells.Item(3, 19).Value = 0.2471144495547739
$ws.Cells.Item(3, 20).Value = 0.247114449554774

# --- Row 4 (new): MuSCs -> Pomc/Oprd1 -> FAPs ---
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.04954
$ws.Cells.Item(4, 8).Value = 0.14862
$ws.Cells.Item(4, 9).Value = 0.08804413699625062
$ws.Cells.Item(4, 10).Value = 0.08804413699625063
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02165966666666666
$ws.Cells.Item(4, 14).Value = 0.064979
$ws.Cells.Item(4, 15).Value = 0.7290280598220596
$ws.Cells.Item(4, 16).Value = 0.7290280598220598
$ws.Cells.Item(4, 17).Value = 0.001073019886666666
$ws.Cells.Item(4, 18).Value = 0.00965717898
$ws.Cells.Item(4, 19).Value = 0.06418664637308422
$ws.Cells.Item(4, 20).Value = 0.06418664637308423

# --- Row 5 (new): MuSCs -> Pomc/Oprd1 -> Resolving-Mac ---
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.04954
$ws.Cells.Item(5, 8).Value = 0.14862
$ws.Cells.Item(5, 9).Value = 0.08804413699625062
$ws.Cells.Item(5, 10).Value = 0.08804413699625063
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.008050666666666666
$ws.Cells.Item(5, 14).Value = 0.024152
$ws.Cells.Item(5, 15).Value = 0.2709719401779404
$ws.Cells.Item(5, 16).Value = 0.2709719401779404
$ws.Cells.Item(5, 17).Value = 0.0003988300266666666
$ws.Cells.Item(5, 18).Value = 0.00358947024
$ws.Cells.Item(5, 19).Value = 0.02385749062316641
$ws.Cells.Item(5, 20).Value = 0.02385749062316641
